$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 (header row): add "Part Number" header in E16, merged across E16:G16 ---
# Copy the look of the existing "Reference Invoice"/"Quantity" header cell (B16) onto E16
$ws.Range("B16").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Value = "Part Number"

# G16 becomes a plain trailing border cell, matching the style used by G10/G11/G13
$ws.Range("G13").Copy()
$ws.Range("G16").PasteSpecial(-4122)

$ws.Range("E16:G16").Merge()

# --- Row 17 (data row): add "{booking:part_number}" value in E17, merged across E17:G17 ---
# Copy the look of the existing "{booking:qty}" data cell (D17) onto E17, then fix alignment to left
$ws.Range("D17").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E17").Value = "{booking:part_number}"
$ws.Range("E17").HorizontalAlignment = -4131

# G17 becomes a plain trailing border cell, matching the style used by G10/G11/G13
$ws.Range("G13").Copy()
$ws.Range("G17").PasteSpecial(-4122)

$ws.Range("E17:G17").Merge()

$excel.CutCopyMode = 0
